# Update the division problems in the table to the new set of values.
# Each cell's text (e.g. "58÷2=") is replaced with a new expression
# (e.g. "29÷4="). Because some of the new values coincide with old
# values used elsewhere in the document, we first replace every "old"
# value with a unique placeholder token, and only then replace every
# placeholder with its final "new" value. This two-phase approach
# avoids any chance of a later replacement accidentally matching text
# that was produced by an earlier replacement.

$d = $word.ActiveDocument

$pairs = @(
    @("58÷2=", "29÷4="),
    @("90÷4=", "35÷4="),
    @("54÷3=", "46÷5="),
    @("10÷9=", "56÷7="),
    @("58÷7=", "62÷4="),
    @("21÷6=", "72÷3="),
    @("33÷4=", "41÷2="),
    @("56÷5=", "13÷4="),
    @("20÷6=", "22÷6="),
    @("80÷8=", "80÷9="),
    @("20÷7=", "99÷2="),
    @("80÷6=", "74÷5="),
    @("88÷2=", "32÷3="),
    @("71÷7=", "22÷3="),
    @("25÷7=", "41÷9="),
    @("43÷9=", "10÷3="),
    @("68÷3=", "52÷5="),
    @("83÷3=", "13÷4="),
    @("10÷7=", "38÷7="),
    @("28÷2=", "50÷7="),
    @("98÷2=", "39÷9="),
    @("18÷7=", "28÷8="),
    @("96÷9=", "70÷2="),
    @("59÷4=", "15÷3="),
    @("69÷4=", "20÷7=")
)

# Phase 1: old -> unique placeholder
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $old = $pairs[$i][0]
    $placeholder = "§PLACEHOLDER_$i§"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $placeholder, 2)
}

# Phase 2: placeholder -> new value
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $new = $pairs[$i][1]
    $placeholder = "§PLACEHOLDER_$i§"
    $d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
